# M0078DDG-336 trimmed white spaces in Shielded Plug diagrams
#
# Re-crop / re-position the two "Shielded Plug" swatch pictures on slide 1
# so the white margin baked into ppt/media/image1.png is trimmed out of the
# visible picture, and shrink the picture frames to match the new (smaller)
# visible area.
#
# NOTE on units: the PowerPoint COM object model works in points (1 pt =
# 12700 EMU) for Left/Top/Width/Height, and in points for
# PictureFormat.CropTop/CropBottom/CropLeft/CropRight (the OOXML
# <a:srcRect> attributes are in 1/1000ths of a percent of the native image
# extent; the source PNG is 256x256 px -> 192x192 pt at 96 dpi, so
# crop_pt = pct_1000ths / 100000 * 192).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Picture 5 (Shapes.Item(4)): top swatch
#   srcRect t="-37477" b="37477"  ->  t="776" b="37477"
#   off    5273436,2149859  ext 1152708,1152708
#                             -> off 5273436,2590799  ext 1152708,711767
# ---------------------------------------------------------------------
$pic5 = $s.Shapes.Item(4)
$pic5.PictureFormat.CropTop    = 1.4899200000000001
$pic5.PictureFormat.CropBottom = 71.95584
$pic5.Left   = 415.2311811023622
$pic5.Top    = 203.99992125984252
$pic5.Width  = 90.7644094488189
$pic5.Height = 56.04464566929134

# ---------------------------------------------------------------------
# Picture 6 (Shapes.Item(5)): bottom swatch
#   srcRect t="64611" b="-64611"  ->  t="64611" b="5282"
#   off    5273436,3302567  ext 1152709,1152709
#                             -> off 5273436,3302568  ext 1152709,347048
# ---------------------------------------------------------------------
$pic6 = $s.Shapes.Item(5)
$pic6.PictureFormat.CropTop    = 124.05311999999999
$pic6.PictureFormat.CropBottom = 10.14144
$pic6.Left   = 415.2311811023622
$pic6.Top    = 260.0447244094488
$pic6.Width  = 90.76448818897637
$pic6.Height = 27.326614173228347

# ---------------------------------------------------------------------
# The cached "today" text on the auto-updating datetimeFigureOut date
# field (Date Placeholder 3) on the slide master and every slide layout
# was re-stamped by PowerPoint (30/03/2020 -> 12/10/2020) the next time
# the deck was saved. Refresh the cached text on each of them.
# ---------------------------------------------------------------------
$newDate = "12/10/2020"

$master = $s.Master
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $layout.Shapes.Item(3).TextFrame.TextRange.Text = $newDate
}
